$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.622.58'
$ws.Range('E2').Value = '  -3.87%  '

$ws.Range('D3').Value = '2.911.32'
$ws.Range('E3').Value = '  -3.54%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.40%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -2.23%  '

$ws.Range('D9').Value = '2.919.68'
$ws.Range('E9').Value = '  -3.51%  '

$ws.Range('E10').Value = '  -5.42%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.86'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.54%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.87%  '

$ws.Range('D13').Value = '3.418.88'
$ws.Range('E13').Value = '  -3.73%  '

$ws.Range('E14').Value = '  +1.07%  '

$ws.Range('D15').Value = '60.740.65'
$ws.Range('E15').Value = '  -3.74%  '

$ws.Range('E16').Value = '  -5.91%  '

$ws.Range('D17').Value = '2.910.60'
$ws.Range('E17').Value = '  -3.69%  '

$ws.Range('E18').Value = '  -6.67%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.35%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.74%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '360.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.92%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.61'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.61%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.90%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.73%  '

$ws.Range('D26').Value = '3.030.22'
$ws.Range('E26').Value = '  -3.86%  '

$ws.Range('E27').Value = '  -3.96%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.177'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.00%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.41%  '

$ws.Range('D30').Value = '0.0₃0856'
$ws.Range('E30').Value = '  -12.53%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.63'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -12.18%  '

$ws.Range('E33').Value = '  -5.12%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.08%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '150.90'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.54%  '

$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.55'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.46%  '

$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.09%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.993'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.63%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.96%  '

$ws.Range('D41').Value = '2.325.58'
$ws.Range('E41').Value = '  -7.79%  '

$ws.Range('E42').Value = '  -8.49%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.84%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.644'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.73%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.57%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0566'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.64%  '

$ws.Range('E48').Value = '  -4.96%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.54%  '

$ws.Range('E50').Value = '  -6.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0921'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.82%  '
